$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly row is inserted before the old row 22, so everything from the
# old row 22 downward shifts one row. Insert a fresh row at 22 and copy the
# old row 24 data (which ends up landing on row 25 after the shift) into the
# new last row, then fix up the dates that moved between rows 22-24.

# Shift rows 22:24 down by one to free row 22 for the new weekly entry, and
# to push the old row 24 down to row 25.
$ws.Rows("22").Insert()

# The freshly inserted blank row 22 pushed the former row22 -> row23,
# former row23 -> row24, former row24 -> row25. Copy the values back down
# into their shifted rows explicitly (values only, since this sheet has no
# per-row formatting beyond the date column style already present).

# Row 25 (previously row 24, now shifted down) - restore its original content.
$ws.Range("A25").Value = 10
$ws.Range("B25").Value = "Vega Modelo de Temuco"
$ws.Range("C25").Value = "La Araucanía"
$ws.Range("D25").Value = 44767
$ws.Range("E25").Value = 9
$ws.Range("F25").Value = 100112036
$ws.Range("G25").Value = "Caigua"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 50
$ws.Range("K25").Value = 20000
$ws.Range("L25").Value = 20000
$ws.Range("M25").Value = 20000
$ws.Range("N25").Value = "$/caja 15 kilos"
$ws.Range("O25").Value = "Región de Arica y Parinacota"
$ws.Range("P25").Value = 1333
$ws.Range("Q25").Value = 15
$ws.Range("R25").Value = "Hortaliza"

# Row 24 (previously row 23, now shifted down) - restore its content, date changes to 44826.
$ws.Range("A24").Value = 10
$ws.Range("B24").Value = "Vega Modelo de Temuco"
$ws.Range("C24").Value = "La Araucanía"
$ws.Range("D24").Value = 44826
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = 100112036
$ws.Range("G24").Value = "Caigua"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 50
$ws.Range("K24").Value = 20000
$ws.Range("L24").Value = 20000
$ws.Range("M24").Value = 20000
$ws.Range("N24").Value = "$/caja 15 kilos"
$ws.Range("O24").Value = "Región de Arica y Parinacota"
$ws.Range("P24").Value = 1333
$ws.Range("Q24").Value = 15
$ws.Range("R24").Value = "Hortaliza"

# Row 23 (previously row 22, now shifted down) - restore its content, date changes to 44755.
$ws.Range("A23").Value = 10
$ws.Range("B23").Value = "Vega Modelo de Temuco"
$ws.Range("C23").Value = "La Araucanía"
$ws.Range("D23").Value = 44755
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = 100112036
$ws.Range("G23").Value = "Caigua"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 50
$ws.Range("K23").Value = 20000
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = 20000
$ws.Range("N23").Value = "$/caja 15 kilos"
$ws.Range("O23").Value = "Región de Arica y Parinacota"
$ws.Range("P23").Value = 1333
$ws.Range("Q23").Value = 15
$ws.Range("R23").Value = "Hortaliza"

# Row 22 - the brand new weekly entry.
$ws.Range("A22").Value = 10
$ws.Range("B22").Value = "Vega Modelo de Temuco"
$ws.Range("C22").Value = "La Araucanía"
$ws.Range("D22").Value = 44841
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = 100112036
$ws.Range("G22").Value = "Caigua"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = 16000
$ws.Range("L22").Value = 16000
$ws.Range("M22").Value = 16000
$ws.Range("N22").Value = "$/caja 15 kilos"
$ws.Range("O22").Value = "Región de Arica y Parinacota"
$ws.Range("P22").Value = 1067
$ws.Range("Q22").Value = 15
$ws.Range("R22").Value = "Hortaliza"

# Apply the date number format (style used by the other date cells in column D)
$ws.Range("D22:D25").NumberFormat = $ws.Range("D21").NumberFormat
